$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Fecha" (D) and "Volumen" (J) values for rows 2 and 3 were swapped.
$d2 = $ws.Range("D2").Value()
$d3 = $ws.Range("D3").Value()
$j2 = $ws.Range("J2").Value()
$j3 = $ws.Range("J3").Value()

$ws.Range("D2").Value = $d3
$ws.Range("D3").Value = $d2

$ws.Range("J2").Value = $j3
$ws.Range("J3").Value = $j2
